$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("DataSet")

# Insert a new row at row 3 (pushes existing row 3 "EmpAccountDetails" and
# everything below it down to row 4+). Excel's default row insert behavior
# copies formatting from the row above (row 2), which matches the target
# styles for the new row exactly.
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the "New Account Details" test case data.
$ws.Range("A3").Value = "New Account Details"
$ws.Range("B3").Value = "avayugundla+3@helenoftroy.com"
$ws.Range("C3").Value = "avayugundla+3@helenoftroy.com"
$ws.Range("D3").Value = "Lotuswave@123"
$ws.Range("E3").Value = "Lotuswave@123"
$ws.Range("F3").Value = "QA"
$ws.Range("G3").Value = "TEST"

# Add mailto/password hyperlinks for the new cells, mirroring how the
# equivalent columns are linked on the other rows of this sheet.
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:avayugundla+3@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:avayugundla+3@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:Lotuswave@123")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:Lotuswave@123")

# Reflect the cursor position the author had active when saving.
$ws.Activate()
$ws.Range("B8").Select()
